# Adds the "citi_cash_balances" field-mapping block to the mapping sheet
# (7 new rows appended after the existing "citi_positions_valuations" rows),
# matching the source commit "added cash balance import logics".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# database_field | report_field | source | table
$newRows = @(
    @("currency",            "CCY",               "Citi", "citi_cash_balances"),
    @("account_id",          "Account ID",         "Citi", "citi_cash_balances"),
    @("account_name",        "Account Name",       "Citi", "citi_cash_balances"),
    @("period ",             "As of Date",         "Citi", "citi_cash_balances"),
    @("opening_balance",     "Opening Balance",     "Citi", "citi_cash_balances"),
    @("ledger_balance",      "Ledger Balance",      "Citi", "citi_cash_balances"),
    @("available_balance",   "Available Balance",   "Citi", "citi_cash_balances")
)

$startRow = 263
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}

$endRow = $startRow + $newRows.Count - 1

# The new block (including the blank column E) picks up the plain-black
# font formatting carried over with the imported data.
$ws.Range("A$($startRow):E$endRow").Font.Color = 0

# Restore the active selection/scroll position to where the new rows were
# added, matching the author's final cursor position.
$ws.Range("C280").Select()
